# Apply the edit described by the diff:
#  - Fill column A for rows 74..94 with sequential numbers 72..92,
#    using the same "row number" cell style already used elsewhere in
#    column A (e.g. A73, which holds 71).
#  - Move the sheet view: drop the scrolled topLeftCell and move the
#    active selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the visual format of the existing numbered cell A73 (style used for
# the sequential index column) so the new A74:A94 values render exactly like
# the rest of the list's numbering column, without introducing a brand new
# style definition.
$ws.Range("A73").Copy()

$startRow = 74
$endRow = 94
$startNumber = 72

for ($r = $startRow; $r -le $endRow; $r++) {
    $target = $ws.Range("A$r")
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $target.Value = $startNumber + ($r - $startRow)
}

$excel.CutCopyMode = 0

# Reset the sheet view: no frozen/scrolled top-left cell, selection on B1.
$ws.Range("B1").Select()
